$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column F (District) so it becomes the new
# "Address" column, and the existing District column shifts to G.
$ws.Columns("F:F").Insert()

# Header
$ws.Range("F2").Value = "Address"

# Address values for rows 3-25 (aligned with the shifted District column in G)
$addresses = @(
    "Govt. High SchoolToranadinniManvi",
    "Sajjalashri High School MegalapetMudagallaLingasagur",
    "Sri Basaveshwar High School Lingasagur",
    "Govt. P U College (High School Section) Hutti Gold MinesLingasugur",
    "Sri Amreswar Girls High School Lingasugur",
    "Govt. Urdu High SchoolLingasugur",
    "Morarji Desai Residential School DevarabhupurLingasagur",
    "G H S Rajalabanda",
    "G H S Alkod",
    "G H S NagarahalLingasagur",
    "G H S ChagabhaviManvi",
    "Morarji Desai Residencial school MaskiLingasugur",
    "Govt. Boys Junior College (High School Section) Sindhanur",
    "G H S Maddipet",
    "G H S R H Colony No – 2Sindhanur",
    "G H S B GanekalDeovadurga",
    "Hamdard High School",
    "Govt. Comp. Jr College for boys High School sectionmaskiLingsugur",
    "G H S HalapurManvi",
    "Kittur Rani Channamma Residensial SchoolAdavibhaviLingasugur",
    "G H S PothnalManvi",
    "Govt. High School BallatagiManvi",
    "Sri Amareshwara High School EachanalLingasur"
)

for ($i = 0; $i -lt $addresses.Length; $i++) {
    $row = 3 + $i
    $ws.Range("F$row").Value = $addresses[$i]
}
